# Add run sheet scheduling V3
# - Insert a new "Worker_v3" worksheet (positioned just before "Room")
#   containing an expanded worker roster (Name / Worker_Id / Efficiency / Wage_Hour).
# - Update the Constraints!Clinic_Close value from 20:00 to 23:00.
# - Update selection/active-sheet UI state to match the edit session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Worker_v3" sheet right before "Room".
# ---------------------------------------------------------------------------
$roomSheet = $wb.Worksheets.Item("Room")
$newSheet = $wb.Worksheets.Add($roomSheet)
$newSheet.Name = "Worker_v3"

# Column A (Name) - write top-to-bottom first so new shared strings are
# interned in row order.
$names = @("Name", "Meddeline", "Marion", "Klarion", "Cindy", "Arlong", "Marlo", "Reno", "Lorean", "Jannik", "Shane", "Mary", "Nashvile", "Larry", "Meghan", "Sharon", "Ello", "Myrion", "Arthur", "Darius")
for ($i = 0; $i -lt $names.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $names[$i]
}

# Column C header (Efficiency).
$newSheet.Cells.Item(1, 3).Value = "Efficiency"

# Column B (Worker_Id).
$ids = @("Worker_Id", "D1", "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 2).Value = $ids[$i]
}

# Column C data (Efficiency values).
$efficiency = @(0, 2, 4, 6, 8, 10, 12, 14, 15, 16, 18, 20, 22, 24, 26, 28, 30, 35, 40)
for ($i = 0; $i -lt $efficiency.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 3).Value = $efficiency[$i]
}

# Column D (Wage_Hour header + values).
$newSheet.Cells.Item(1, 4).Value = "Wage_Hour"
$wage = @(35000, 38000, 45000, 45000, 48000, 48000, 49000, 49000, 50000, 56000, 62000, 65000, 65000, 68000, 70000, 75000, 80000, 85000, 90000)
for ($i = 0; $i -lt $wage.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 4).Value = $wage[$i]
}

$newSheet.Range("C24").Select()

# ---------------------------------------------------------------------------
# 2. Update the Clinic_Close constraint (20:00 -> 23:00).
# ---------------------------------------------------------------------------
$wsConstraints = $wb.Worksheets.Item("Constraints")
$wsConstraints.Range("B14").Value = 23 / 24

# ---------------------------------------------------------------------------
# 3. Restore selection / active-sheet UI state.
# ---------------------------------------------------------------------------
$wsWorker = $wb.Worksheets.Item("Worker")
$wsWorker.Range("L23").Select()

$wsConstraints.Activate()
$wsConstraints.Range("B3").Select()
